$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text (non-numeric-looking) columns: set directly
$textUpdates = @{
    "B10" = "MandalaExchangeToken"
    "C10" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
    "B11" = "BitrueCoin"
    "C11" = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
    "B12" = "BitMartToken"
    "C12" = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
    "B13" = "BitForexToken"
    "C13" = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
    "B14" = "One"
    "C14" = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
    "B15" = "TigerCash"
    "C15" = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
    "B16" = "LEO"
    "C16" = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
    "B17" = "GateToken"
    "C17" = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
    "B18" = "BTSEToken"
    "C18" = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
    "B19" = "BitpandaEcosystemToken"
    "C19" = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
    "B20" = "LiechtensteinCryptoassetsExchange"
    "C20" = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
    "B41" = "KickToken"
    "C41" = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
    "B42" = "BKEXToken"
    "C42" = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
}

# Numeric-looking columns that must stay text: force Text format, set, then clear style
$numericTextUpdates = @{
    "D2" = "257.73"
    "E2" = "5.14%"
    "G2" = "20"
    "D3" = "27.66"
    "E3" = "-2.81%"
    "G3" = "20"
    "D4" = "5.232"
    "E4" = "-0.25%"
    "G4" = "20"
    "D5" = "0.05917"
    "E5" = "3.80%"
    "G5" = "20"
    "D6" = "6.688"
    "E6" = "1.15%"
    "G6" = "20"
    "D7" = "0.8685"
    "E7" = "2.20%"
    "G7" = "20"
    "D8" = "1.053"
    "E8" = "17.77%"
    "G8" = "20"
    "D9" = "0.1419"
    "E9" = "3.50%"
    "G9" = "20"
    "D10" = "0.07194"
    "E10" = "2.26%"
    "G10" = "20"
    "D11" = "0.03257"
    "E11" = "3.09%"
    "G11" = "20"
    "D12" = "0.09223"
    "E12" = "0.19%"
    "G12" = "20"
    "D13" = "0.001545"
    "E13" = "1.45%"
    "G13" = "20"
    "D14" = "0.0006065"
    "E14" = "2.01%"
    "G14" = "20"
    "D15" = "0.005860"
    "E15" = "-2.31%"
    "G15" = "20"
    "D16" = "3.485"
    "E16" = "-0.18%"
    "G16" = "20"
    "D17" = "3.270"
    "E17" = "2.19%"
    "G17" = "20"
    "D18" = "2.227"
    "E18" = "2.41%"
    "G18" = "20"
    "D19" = "0.3151"
    "E19" = "-0.64%"
    "G19" = "20"
    "D20" = "0.03614"
    "E20" = "9.85%"
    "G20" = "20"
    "D21" = "0.1291"
    "E21" = "0.28%"
    "G21" = "20"
    "D22" = "3.530"
    "E22" = "0.77%"
    "G22" = "20"
    "D23" = "0.04180"
    "E23" = "2.20%"
    "G23" = "20"
    "D24" = "0.1400"
    "E24" = "1.57%"
    "G24" = "20"
    "D25" = "0.001219"
    "E25" = "-0.03%"
    "G25" = "20"
    "D26" = "0.004532"
    "E26" = "9.20%"
    "G26" = "20"
    "D27" = "0.0001201"
    "E27" = "0.13%"
    "G27" = "20"
    "D28" = "0.0001940"
    "E28" = "33.99%"
    "G28" = "20"
    "G29" = "20"
    "G30" = "20"
    "G31" = "20"
    "G32" = "20"
    "G33" = "20"
    "G34" = "20"
    "G35" = "20"
    "G36" = "20"
    "G37" = "20"
    "G38" = "20"
    "G39" = "20"
    "D40" = "0.03818"
    "G40" = "20"
    "D41" = "0.005484"
    "E41" = "6.35%"
    "G41" = "20"
    "D42" = "0.1107"
    "E42" = "4.01%"
    "G42" = "20"
    "D43" = "0.002382"
    "E43" = "8.33%"
    "G43" = "20"
    "D44" = "0.009888"
    "E44" = "7.99%"
    "G44" = "20"
    "D45" = "0.00005437"
    "E45" = "2.92%"
    "G45" = "20"
    "E46" = "0.14%"
    "G46" = "20"
    "E47" = "4.08%"
    "G47" = "20"
    "D48" = "0.002161"
    "E48" = "-4.74%"
    "G48" = "20"
    "D49" = "0.00002102"
    "E49" = "0.14%"
    "G49" = "20"
    "D50" = "0.0002002"
    "E50" = "0.14%"
    "G50" = "20"
    "G51" = "20"
}

foreach ($addr in $textUpdates.Keys) {
    $ws.Range($addr).Value = $textUpdates[$addr]
}

foreach ($addr in $numericTextUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $numericTextUpdates[$addr]
    $cell.Style = "Normal"
}
